# Week 5 Meeting.pptx — applies:
#  1) Update the cached "last modified" date text shown in the date/time
#     placeholder field on every slide layout + the slide master from
#     10/31/2021 -> 12/11/2021 (mirrors the footer date fields in the diff).
#  2) Split the subtitle run "Week 5 Meeting, 1/11/2021" on slide 1 into two
#     runs: "Week 5 " + "Meeting, 1/11/2021" (same visible text, matches the
#     run-split in the diff).

$p = $ppt.ActivePresentation

$oldDate = "10/31/2021"
$newDate = "12/11/2021"

function Update-DateShape($shp) {
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# -- Slide master date placeholder --
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# -- Every slide layout's date placeholder --
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# -- Slide 1 subtitle: split "Week 5 Meeting, 1/11/2021" into two runs --
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text.StartsWith("Week 5 Meeting, 1/11/2021")) {
            $firstPart = $tr.Characters(1, 7)
            $firstPart.Text = "Week 5 "
        }
    }
}
